$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58: B58 161 -> 162, D58 recalculated (C58 unchanged at 244)
$ws.Range("B58").Value = 162
$ws.Range("D58").Value = 66.39344262295081

# Row 59: B59 1 -> 11, C59 227 -> 228, D59 recalculated
$ws.Range("B59").Value = 11
$ws.Range("C59").Value = 228
$ws.Range("D59").Value = 4.824561403508771
